# Insert a new weekly data row before the current row 198, shifting the
# existing rows 198-247 down to 199-248 (dimension grows to A1:R248).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(198).Insert()

$ws.Cells.Item(198, 1).Value = 7
$ws.Cells.Item(198, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(198, 3).Value = "Ñuble"
$ws.Cells.Item(198, 4).Value = 44782
$ws.Cells.Item(198, 5).Value = 16
$ws.Cells.Item(198, 6).Value = 100112043
$ws.Cells.Item(198, 7).Value = "Pepino ensalada"
$ws.Cells.Item(198, 8).Value = "Sin especificar"
$ws.Cells.Item(198, 9).Value = "Primera"
$ws.Cells.Item(198, 10).Value = 120
$ws.Cells.Item(198, 11).Value = 20000
$ws.Cells.Item(198, 12).Value = 21000
$ws.Cells.Item(198, 13).Value = 20500
$ws.Cells.Item(198, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(198, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(198, 16).Value = 342
$ws.Cells.Item(198, 17).Value = 60
$ws.Cells.Item(198, 18).Value = "Hortaliza"
